# Update the metrics_6_6 sheet: columns B..Q for rows 2..26 all share the
# same new metric values (the whole block was re-computed and every row
# got the same new figures, mirroring the original data which also had
# identical values across all rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.5692387035740197,    # B - r2
    0.2355707868591003,    # C - r2_sup
    0.6363693091229717,    # D - r2_test
    -1.465209168673161,    # E - r2_val
    0.3088482927627829,    # F - r2_vt
    0.2557185134775685,    # G - mse
    0.4537982025430372,    # H - mse_sup
    0.1961375308550407,    # I - mse_test
    0.4876673787606747,    # J - mse_val
    0.3419024548078576,    # K - mse_vt
    0.2823828768330753,    # L - mape
    0.5056861808251918,    # M - rmse
    0.06015717143422494,   # N - r2_adj
    0.5272142884576241,    # O - rsd
    28.72735599237364,     # P - aic
    44.57274171566024      # Q - bic
)

$firstRow = 2
$lastRow = 26
$firstCol = 2   # column B
$lastCol = 17   # column Q

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - $firstCol]
    }
}
